# Add a new "AuthData" worksheet (placed after the existing "UserData" sheet)
# holding the data for an OAuth client-credentials authentication test, used
# to build a custom REST-assured request specification.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "AuthData"
$ws.Activate()
$excel.ActiveWindow.DisplayGridlines = $false

# ---- values --------------------------------------------------------------
# Force text number format on every data cell first so the numeric-looking
# value ("1711") is stored as text, matching the other text columns.
$ws.Range("A1:C2").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "testname"
$ws.Range("B1").Value = "formParam"
$ws.Range("C1").Value = "user_id"

# Data row
$ws.Range("A2").Value = "testAuthWithClientCredentials"
$ws.Range("B2").Value = "client_id:rest_assured_oauth_demo_app;client_secret:05c3e7b288e0f9cd2a51afa1dc660d11;grant_type:client_credentials"
$ws.Range("C2").Value = "1711"

# ---- row heights / column widths (match the rest of the workbook) --------
$ws.Rows.Item(1).RowHeight = 13.65
$ws.Rows.Item(2).RowHeight = 13.65
$ws.Range("A:C").ColumnWidth = 15.5

# ---- cell colouring / borders ---------------------------------------------
$highlight = 0xDBDBDB
$white = 0xFFFFFF
$gray = 0xAAAAAA
$lightGray = 0xA5A5A5
$darkGray = 0x3F3F3F

function Set-Box($rng, $left, $top, $right, $bottom) {
    $rng.Borders.Item(7).LineStyle = 1
    $rng.Borders.Item(7).Weight = 2
    $rng.Borders.Item(7).Color = $left

    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(8).Color = $top

    $rng.Borders.Item(10).LineStyle = 1
    $rng.Borders.Item(10).Weight = 2
    $rng.Borders.Item(10).Color = $right

    if ($bottom -eq $null) {
        $rng.Borders.Item(9).LineStyle = 0
    } else {
        $rng.Borders.Item(9).LineStyle = 1
        $rng.Borders.Item(9).Weight = 2
        $rng.Borders.Item(9).Color = $bottom
    }
}

# Header row: boxed cells, first column highlighted
$ws.Range("A1").Interior.Color = $highlight
Set-Box $ws.Range("A1") $gray $gray $gray $gray

$ws.Range("B1").Interior.Color = $white
Set-Box $ws.Range("B1") $gray $gray $gray $gray

$ws.Range("C1").Interior.Color = $white
Set-Box $ws.Range("C1") $gray $gray $gray $null

# Data row
$ws.Range("A2").Interior.Color = $highlight
Set-Box $ws.Range("A2") $lightGray $gray $darkGray $lightGray

$ws.Range("B2").Interior.Color = $white
Set-Box $ws.Range("B2") $darkGray $gray $lightGray $lightGray

$ws.Range("C2").Interior.Color = $white
$ws.Range("C2").Borders.Item(7).LineStyle = 1
$ws.Range("C2").Borders.Item(7).Weight = 2
$ws.Range("C2").Borders.Item(7).Color = $lightGray
$ws.Range("C2").Borders.Item(8).LineStyle = 0
$ws.Range("C2").Borders.Item(10).LineStyle = 1
$ws.Range("C2").Borders.Item(10).Weight = 2
$ws.Range("C2").Borders.Item(10).Color = $gray
$ws.Range("C2").Borders.Item(9).LineStyle = 1
$ws.Range("C2").Borders.Item(9).Weight = 2
$ws.Range("C2").Borders.Item(9).Color = $gray

# ---- sheet / page setup (match the rest of the workbook) -----------------
$ws.PageSetup.LeftMargin = 72
$ws.PageSetup.RightMargin = 72
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 18
$ws.PageSetup.FooterMargin = 18
$ws.PageSetup.Zoom = 100
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.Orientation = 1
$ws.PageSetup.CenterFooter = "&`"Helvetica Neue,Regular`"&12&K000000&P"
